# Adds the four new "3V Model" content slides (Volume/Velocity/Variety breakdown)
# to the end of the deck, each using the "Title and Content" layout (index 2),
# matching the other body slides already in the presentation.

$p = $ppt.ActivePresentation

$newSlidesContent = @(
    @{ Title = '3V Model'; Body = @(
        'Volume',
        'Velocity',
        'Variety'
    ); Autofit = $false },
    @{ Title = 'Volume'; Body = @(
        'Data Volume ',
        '44x increase from 2009 2020 ',
        'From 0.8 zettabytes to 35zb ',
        'Data volume is increasing exponentially  day by day',
        'By 2020, International Data Corporation predicts the number will reach 40,000 EB, or 40 Zettabytes (ZB) . ',
        'The world’s information is doubling every two years. By 2020, there will be 5,200 GB of data for every person on Earth. ',
        'By 2020, the amount of high-value data worth analyzing will double and 60% of information delivered to decision makers will be actionable.'
    ); Autofit = $true },
    @{ Title = 'Variety'; Body = @(
        'XML',
        'JSON',
        'CSV ',
        'TEXT',
        'Parquet',
        'AVRO',
        'Relational Database',
        'Non Relational Database'
    ); Autofit = $false },
    @{ Title = 'Variety'; Body = @(
        'Streaming data – the data which got ability to change randomly per sec/Per minute based',
        'Normal data - > custom pull data, whenever we want data, we can manually pull it and use it analytical purpose'
    ); Autofit = $false }
)

foreach ($slideInfo in $newSlidesContent) {
    # ppLayoutText (2) = "Title and Content" layout, appended at the end of the deck
    $newSlide = $p.Slides.Add($p.Slides.Count + 1, 2)

    $titleShape = $newSlide.Shapes.Item(1)
    $titleShape.TextFrame.TextRange.Text = $slideInfo.Title

    $bodyShape = $newSlide.Shapes.Item(2)
    $bodyText = [string]::Join("`r", $slideInfo.Body)
    $bodyShape.TextFrame.TextRange.Text = $bodyText

    if ($slideInfo.Autofit) {
        $bodyShape.TextFrame.AutoSize = 2
    }
}

